$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C3:F3 hold numeric-looking text ("29", "0.0"); force the Text number
# format on just those cells so Excel keeps them as strings instead of
# silently coercing them to numbers. A3/B3 are plain text already and
# keep the sheet's default style.
$ws.Range("C3:F3").NumberFormat = "@"

$ws.Range("A3").Value = "2024-09-25T17:53:44Z"
$ws.Range("B3").Value = "temperature"
$ws.Range("C3").Value = "29"
$ws.Range("D3").Value = "0.0"
$ws.Range("E3").Value = "0.0"
$ws.Range("F3").Value = "0.0"
